$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> epochs, start_run, end_run, best_val_loss
$data = @(
    @{ Row = 2; Epochs = 100; Start = "05/22/2024 00:51:42"; End = "05/22/2024 02:42:14"; Loss = 0.3866017184072858 },
    @{ Row = 3; Epochs = 100; Start = "05/22/2024 02:42:15"; End = "05/22/2024 04:34:42"; Loss = 0.3744101685495928 },
    @{ Row = 4; Epochs = 100; Start = "05/22/2024 04:34:42"; End = "05/22/2024 06:44:40"; Loss = 0.3931281571844637 },
    @{ Row = 5; Epochs = 100; Start = "05/22/2024 06:44:40"; End = "05/22/2024 08:19:29"; Loss = 0.4401574363816773 },
    @{ Row = 6; Epochs = 100; Start = "05/22/2024 08:19:30"; End = "05/22/2024 10:11:48"; Loss = 0.4253510647461318 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 9).Value = $item.Epochs        # column I
    $ws.Cells.Item($r, 32).Value = $item.Start         # column AF
    $ws.Cells.Item($r, 33).Value = $item.End           # column AG
    $ws.Cells.Item($r, 34).Value = $item.Loss          # column AH
}
